$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 20

$ws.Cells.Item(19, 1).Copy($ws.Cells.Item($row, 1)) | Out-Null
$ws.Cells.Item($row, 1).Value = (Get-Date -Year 2025 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0)

$ws.Cells.Item(19, 2).Copy($ws.Cells.Item($row, 2)) | Out-Null
$ws.Cells.Item($row, 2).Value = 0.83199074074074075

$ws.Cells.Item($row, 3).Value = "Test007"
$ws.Cells.Item($row, 4).Value = "Outer"
$ws.Cells.Item($row, 5).Value = 25
$ws.Cells.Item($row, 6).Value = 50
$ws.Cells.Item($row, 7).Value = 1250
$ws.Cells.Item($row, 8).Value = 12546
$ws.Cells.Item($row, 9).Value = 6
$ws.Cells.Item($row, 10).Value = "Cash"
